$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Sprint 1 / TEODORO block
# ---------------------------------------------------------------------------

# Paragraph 2: "TEODORO" + " - Frontend Lead:" -> merge into a single run
# "TEODORO - Frontend Lead:" (still bold) and strike the whole line through.
$pTeo = $d.Paragraphs(2)
$teoRange = $pTeo.Range
$teoStart = $teoRange.Start
$teoEnd = $teoRange.End - 1   # exclude the paragraph mark
$nameLen = ("TEODORO").Length
$nameEnd = $teoStart + $nameLen
# Remove everything after "TEODORO" (i.e. the " - Frontend Lead:" run) …
$tail = $d.Range($nameEnd, $teoEnd)
$tail.Delete()
# … then re-append the same text onto the (now sole) "TEODORO" run so the
# whole paragraph collapses to a single run.
$nameRun = $d.Range($teoStart, $nameEnd)
$nameRun.InsertAfter(" - Frontend Lead:")

# Strike through the whole (merged) paragraph, including the paragraph mark.
$d.Paragraphs(2).Range.Font.StrikeThrough = 1

# The five bullet points under TEODORO also get struck through.
for ($i = 3; $i -le 7; $i++) {
    $d.Paragraphs($i).Range.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# Sprint 1 / ULISES block
# ---------------------------------------------------------------------------

# "Set up Shopify store backend (...)" -> struck through.
$d.Paragraphs(9).Range.Font.StrikeThrough = 1

# "Configure Shopify admin and create documentation for client" -> struck
# through, with "client" split into its own run wrapped in proofErr markers
# (as Word does when it flags a grammar issue on that word).
$pConfig = $d.Paragraphs(10)
$cFull = $pConfig.Range
$cStart = $cFull.Start
$cEnd = $cFull.End - 1
$cRange = $d.Range($cStart, $cEnd)
$innerXml = '<w:r><w:t xml:space="preserve">Configure Shopify admin and create documentation for </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>client</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>'
$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cRange.InsertXML($packageXml)
$d.Paragraphs(10).Range.Font.StrikeThrough = 1

# "Prepare sample product data and categories" -> struck through.
$d.Paragraphs(12).Range.Font.StrikeThrough = 1

# ("Set up authentication flow ..." and "Create initial deployment pipeline
# ..." are left untouched.)

# ---------------------------------------------------------------------------
# Sprint 2 block
# ---------------------------------------------------------------------------

# "Build Landing Page with hero section and ..." -> struck through.
$d.Paragraphs(18).Range.Font.StrikeThrough = 1

# "Create Product Listing page with filtering/sorting" -> only the
# "Create Product Listing page" portion is struck through; the rest
# (" with filtering/sorting") keeps its normal formatting.
$pList = $d.Paragraphs(19)
$lStart = $pList.Range.Start
$lSplitLen = ("Create Product Listing page").Length
$lFirst = $d.Range($lStart, $lStart + $lSplitLen)
$lFirst.Font.StrikeThrough = 1

# "Implement Product Detail page (...)" -> struck through.
$d.Paragraphs(20).Range.Font.StrikeThrough = 1
